$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9749083518981934
$ws.Range("B1").Value = 5.505975246429443
$ws.Range("C1").Value = 4.182441234588623
$ws.Range("D1").Value = 1.03018856048584
$ws.Range("E1").Value = 0.6485534310340881
